$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric-looking identifiers stay stored as text (shared strings),
# matching the original workbook's formatting for these fields.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"

$ws.Range("A2").Value = "EMPAQUES BELEN S.A."
$ws.Range("B2").Value = "3101135332"
$ws.Range("C2").Value = "EMPAQUES BELEN S.A."
$ws.Range("D2").Value = "24385119"
$ws.Range("E2").Value = "pedidos@empaquesbelen.com"
